# Add five new books (rows 33-37) to the worksheet, matching the author's
# commit "fix: retrieve image links from Google Books API".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Description text blocks (kept as here-strings so embedded quotes/newlines
# are preserved exactly).

$disneyWarDesc = @"
"When You Wish Upon a Star," "Whistle While You Work," "The Happiest Place on Earth" -- these are lyrics indelibly linked to Disney, one of the most admired and best-known companies in the world. So when Roy Disney, chairman of Walt Disney Animation and nephew of founder Walt Disney, abruptly resigned in November 2003 and declared war on chairman and chief executive Michael Eisner, he sent shock waves through the entertainment industry, corporate boardrooms, theme parks, and living rooms around the world -- everywhere Disney does business and its products are cherished.
"DisneyWar" is the breathtaking, dramatic inside story of what drove America's best-known entertainment company to civil war, told by one of our most acclaimed writers and reporters.
Drawing on unprecedented access to both Eisner and Roy Disney, current and former Disney executives and board members, as well as thousands of pages of never-before-seen letters, memos, transcripts, and other documents, James B. Stewart gets to the bottom of mysteries that have enveloped Disney for years: What really caused the rupture with studio chairman Jeffrey Katzenberg, a man who once regarded Eisner as a father but who became his fiercest rival? How could Eisner have so misjudged Michael Ovitz, a man who was not only "the most powerful man in Hollywood" but also his friend, whom he appointed as Disney president and immediately wanted to fire? What caused the break between Eisner and Pixar chairman Steve Jobs, and why did Pixar abruptly abandon its partnership with Disney? Why did Eisner so mistrust Roy Disney that he assigned Disney company executives to spy on him? How did Eisner control the Disney board for so long, and whatreally happened in the fateful board meeting in September 2004, when Eisner played his last cards?
Here, too, is the creative process that lies at the heart of Disney -- from the making of "The Lion King" to "Pirates of the Caribbean." Even as the executive suite has been engulfed in turmoil, Disney has worked -- and sometimes clashed -- with a glittering array of stars, directors, designers, artists, and producers, many of whom tell their stories here for the first time.
Stewart describes how Eisner lost his chairmanship and why he felt obliged to resign as CEO, effective 2006. No other book so thoroughly penetrates the secretive world of the corporate boardroom. "DisneyWar" is an enthralling tale of one of America's most powerful media and entertainment companies, the people who control it, and those trying to overthrow them.
"DisneyWar" is an epic achievement. It tells a story that -- in its sudden twists, vivid, larger-than-life characters, and thrilling climax -- might itself have been the subject of a Disney animated classic -- except that it's all true.
"@

$defenseDesc = @"
Eddie Flynn used to be a con artist. Then he became a lawyer. Turns out the two jobs aren’t all that different. . .
He vowed never to set foot in a courtroom again after a case gone disastrously wrong. But today Eddie doesn’t have a choice. Because this time, it’s personal.
The head of the Russian mob in New York City, on trial for murder, has kidnapped Eddie’s beloved ten-year-old daughter. Now Eddie has exactly forty-eight hours to draw upon his razor-sharp instincts and use every con, bluff, grift, and trick in the book to defend an impossible trial and save his daughter—or die trying. . .
"@

$badBloodDesc = @"
The gripping story of Elizabeth Holmes and Theranos—one of the biggest corporate frauds in history—a tale of ambition and hubris set amid the bold promises of Silicon Valley, rigorously reported by the prize-winning journalist. With a new Afterword covering her trial and sentencing, bringing the story to a close.
In 2014, Theranos founder and CEO Elizabeth Holmes was widely seen as the next Steve Jobs: a brilliant Stanford dropout whose startup “unicorn” promised to revolutionize the medical industry with its breakthrough device, which performed the whole range of laboratory tests from a single drop of blood. Backed by investors such as Larry Ellison and Tim Draper, Theranos sold shares in a fundraising round that valued the company at more than `$9 billion, putting Holmes’s worth at an estimated `$4.5 billion. There was just one problem: The technology didn’t work. Erroneous results put patients in danger, leading to misdiagnoses and unnecessary treatments. All the while, Holmes and her partner, Sunny Balwani, worked to silence anyone who voiced misgivings—from journalists to their own employees.
"@

$thursdayDesc = @"
In a peaceful retirement village, four unlikely friends meet weekly in the Jigsaw Room to discuss unsolved crimes; together they call themselves the Thursday Murder Club.
When a local developer is found dead with a mysterious photograph left next to the body, the Thursday Murder Club suddenly find themselves in the middle of their first live case.
As the bodies begin to pile up, can our unorthodox but brilliant gang catch the killer, before it's too late?
"@

$lordEdgwareDesc = @"
In this official authorized edition from the Queen of Mystery, detective Hercule Poirot must solve a most confounding conundrum when Lord Edgware has a most unnatural death.
When Lord Edgware is found murdered the police are baffled. His estranged actress wife was seen visiting him just before his death and Hercule Poirot himself heard her brag of her plan to "get rid" of him.
But how could she have stabbed Lord Edgware in his library at exactly the same time she was seen dining with friends? It's a case that almost proves to be too much for the great Poirot.
"@

# New row data: id, title, author, description, imageurl, date(serial)
$rows = @(
  @{ Row=33; Id=32; Title="Disney War"; Author="James B. Stewart"; Desc=$disneyWarDesc; Image="32 Disney War.jpg"; Date=45575.967919884264 },
  @{ Row=34; Id=33; Title="The Defense"; Author="Steve Cavanagh"; Desc=$defenseDesc; Image="33 The Defense.jpg"; Date=45578.26440537037 },
  @{ Row=35; Id=34; Title="Bad Blood: Secrets and Lies in a Silicon Valley Startup"; Author="John Carreyrou"; Desc=$badBloodDesc; Image="34 Bad Blood.jpg"; Date=45584.20642707176 },
  @{ Row=36; Id=35; Title="The Thursday Murder Club"; Author="Richard Osman"; Desc=$thursdayDesc; Image="35 The Thursday Murder Club.jpg"; Date=45586.123904548615 },
  @{ Row=37; Id=36; Title="Lord Edgware Dies"; Author="Agatha Christie"; Desc=$lordEdgwareDesc; Image="36 Lord Edgware Dies.jpg"; Date=45589.084046307871 }
)

foreach ($r in $rows) {
  # Clone formatting (number formats / styles) from the previous row first,
  # then overwrite values so new cells pick up the same styles used by the
  # rest of the table (integer id column, date column, etc.).
  $srcRow = $r.Row - 1
  $ws.Range("A" + $srcRow + ":F" + $srcRow).Copy()
  $ws.Range("A" + $r.Row + ":F" + $r.Row).PasteSpecial(-4122)

  $ws.Range("A" + $r.Row).Value = $r.Id
  $ws.Range("B" + $r.Row).Value = $r.Title
  $ws.Range("C" + $r.Row).Value = $r.Author
  $ws.Range("D" + $r.Row).Value = $r.Desc
  $ws.Range("E" + $r.Row).Value = $r.Image
  $ws.Range("F" + $r.Row).Value = $r.Date

  # Multi-line descriptions otherwise leave the engine's auto row-height
  # calculation in place; AutoFit restores the sheet's default height so
  # the row doesn't end up with a stray explicit ht/customHeight attribute.
  $ws.Rows($r.Row).AutoFit()
}

$excel.CutCopyMode = 0

# Update the selection so it mirrors the appended rows (whole row 38,
# matching the author's post-edit cursor position).
$ws.Rows(38).Select()
